# Update calendar feed with 3rd entry
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4: TEST/EVENT course entry
$ws.Cells.Item(4, 1).Value = "TEST/EVENT"
$ws.Cells.Item(4, 2).Value = "MYMY"
$ws.Cells.Item(4, 3).Value = "Hello"
$ws.Cells.Item(4, 4).Value = "Tutorial"

# Start Date (E4) - same date style as the existing rows
$ws.Cells.Item(4, 5).Value = 45884
$ws.Cells.Item(4, 5).NumberFormat = $ws.Cells.Item(2, 5).NumberFormat

# Start Time (F4) left blank, but stamped with an h:mm time format
$ws.Cells.Item(4, 6).NumberFormat = "h:mm"

# End Date (G4) - same date style as the existing rows
$ws.Cells.Item(4, 7).Value = 45885
$ws.Cells.Item(4, 7).NumberFormat = $ws.Cells.Item(2, 7).NumberFormat

# Timezone / Location (I4, J4) and Transparent flag (M4)
$ws.Cells.Item(4, 9).Value = "Sydney, Australia"
$ws.Cells.Item(4, 10).Value = "Sydney, Australia"
$ws.Cells.Item(4, 13).Value = "Transparent"

# Move the active selection to I5, matching the saved view state
[void]$ws.Range("I5").Select()
